$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update overall progress percentage (C6): 25% -> 52% ---
# Use a leading apostrophe so Excel keeps it as literal text "52%"
# instead of auto-converting it to a numeric percentage value.
$ws.Range("C6").Value = "'52%"

# --- Mark several tasks as Complete ---
# These rows currently use the "Not Started" look (blue fill) on columns
# A:F, and need to switch to the "Complete" look (green fill), with column
# E's text updated to "Complete". Column G is left untouched, matching the
# source diff. Row 41 already has the desired "Complete" formatting, so we
# copy its A:F formatting onto each target row rather than constructing a
# brand new style.
$completedRows = 27,32,33,34,35,36,37,38,44,45,50

$donor = $ws.Range("A41:F41")
$donor.Copy()

foreach ($r in $completedRows) {
    $target = $ws.Range("A" + $r + ":F" + $r)
    $target.PasteSpecial(-4122)
    $ws.Range("E" + $r).Value = "Complete"
}

$excel.CutCopyMode = 0
